# Apply updated cryptocurrency price/volume figures to Sheet1.
# Cells that hold numeric-looking text (e.g. "596.51") are written with a
# temporary Text number format so Excel keeps them as strings instead of
# silently parsing them into numbers (which would drop formatting such as
# trailing zeros, e.g. "8.10" -> 8.1). The style is reset to "Normal"
# immediately afterwards so no residual formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.438.54"

# Row 3
$ws.Range("D3").Value = "3.197.39"
$ws.Range("E3").Value = "  +1.67%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.75%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.07%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").Value = "3.194.79"
$ws.Range("E8").Value = "  +1.60%  "

# Row 9
$ws.Range("E9").Value = "  +4.29%  "

# Row 10
$ws.Range("E10").Value = "  +1.80%  "

# Row 11
$ws.Range("E11").Value = "  -1.28%  "

# Row 12
$ws.Range("E12").Value = "  +4.14%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000268"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.83%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.92%  "

# Row 15
$ws.Range("D15").Value = "3.720.99"
$ws.Range("E15").Value = "  +1.63%  "

# Row 16
$ws.Range("D16").Value = "66.436.34"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.47%  "

# Row 18
$ws.Range("D18").Value = "3.196.23"
$ws.Range("E18").Value = "  +1.58%  "

# Row 19
$ws.Range("E19").Value = "  +0.98%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "515.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.97%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.53%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.742"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.46%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.44%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.59%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.53%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.36%  "

# Row 27
$ws.Range("E27").Value = "  +5.68%  "

# Row 28
$ws.Range("E28").Value = "  +4.14%  "

# Row 29
$ws.Range("E29").Value = "  +7.82%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +17.75%  "

# Row 31
$ws.Range("E31").Value = "  +3.69%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.33%  "

# Row 33
$ws.Range("E33").Value = "  +3.13%  "

# Row 34
$ws.Range("E34").Value = "  +0.01%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.78%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "513.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.66%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.65%  "

# Row 38
$ws.Range("E38").Value = "  +1.62%  "

# Row 39
$ws.Range("E39").Value = "  +2.88%  "

# Row 40
$ws.Range("E40").Value = "  +3.43%  "

# Row 41
$ws.Range("E41").Value = "  +7.21%  "

# Row 42
$ws.Range("E42").Value = "  -2.06%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.304"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.38%  "

# Row 44
$ws.Range("D44").Value = "0.0₃0674"
$ws.Range("E44").Value = "  +16.88%  "

# Row 45
$ws.Range("E45").Value = "  +2.81%  "

# Row 46
$ws.Range("D46").Value = "2.920.19"
$ws.Range("E46").Value = "  -2.77%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.79%  "

# Row 48
$ws.Range("E48").Value = "  +3.10%  "

# Row 49
$ws.Range("E49").Value = "  +0.01%  "

# Row 50
$ws.Range("E50").Value = "  +5.61%  "

# Row 51
$ws.Range("E51").Value = "  +10.55%  "
